$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "22-03-2025"
$ws.Range("B5").Value = "Kolkata Knight Riders vs Royal Challengers Bengaluru"
$ws.Range("C5").Value = "Kolkata Knight Riders"
$ws.Range("D5").Value = "Kolkata Knight Riders"
$ws.Range("E5").Value = "vijay"
